# edit.ps1 - apply the "Quantum Computing" -> "The Art of Healing" rewrite
# plus the TimesNewToman -> Times New Roman font-name fix.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Global font-name fix: TimesNewToman -> Times New Roman
#    Applying Font.Name on an explicit-bounds Range (rather than the
#    $d.Content object directly) keeps every existing run boundary
#    intact (no merging), it just rewrites <w:rFonts>.
# ---------------------------------------------------------------------
$fullRange = $d.Range(0, $d.Content.End)
$fullRange.Font.Name = "Times New Roman"

# ---------------------------------------------------------------------
# Helper: a tightly scoped Find/Replace so we touch as little as
# possible of the surrounding run structure.
# ---------------------------------------------------------------------
function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Text = ""
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $old"
    }
}

# ---------------------------------------------------------------------
# 2) Title
# ---------------------------------------------------------------------
Replace-Text "Quantum Computing: Unraveling the Enigma" "The Art of Healing: A Journey Through the History of Medicine"

# ---------------------------------------------------------------------
# 3) Author name (byline)
# ---------------------------------------------------------------------
Replace-Text "Eleanor Lawson" "Samantha Adams"

# ---------------------------------------------------------------------
# 4) Email address pieces
# ---------------------------------------------------------------------
Replace-Text "eleanor" "samantha"
Replace-Text "lawson@quantum-computing-lab" "adams@school"
Replace-Text "org" "edu"

# ---------------------------------------------------------------------
# 5) Body paragraph - sentence-by-sentence rewrite (existing runs)
# ---------------------------------------------------------------------
Replace-Text "Step into the realm of quantum computing, where the fabric of reality intertwines with the enigmatic dance of subatomic particles" "From ancient herbal remedies to cutting-edge biotechnology, medicine has a rich and ever-evolving history"

Replace-Text " A universe of mind-bending possibilities unfolds before us, challenging our understanding of computation and promising transformative breakthroughs across diverse disciplines" " It is a tale of human ingenuity, resilience, and the pursuit of understanding our bodies and alleviating suffering"

Replace-Text " Embark on a journey into this fascinating realm, where the quantum world whispers secrets of untapped potential" " Medicine's journey has been shaped by countless individuals--doctors, scientists, healers, and innovators--who dedicated their lives to pushing the boundaries of knowledge and making a difference in people's lives"

Replace-Text "In the heart of a quantum computer, qubits, the quantum counterparts of classical bits, reside in a superposition of states, unlocking the parallel processing of intricate algorithms" "Medicine's origins can be traced back to prehistoric times, with early humans using plants, animal products, and other natural materials to treat illnesses and injuries"

Replace-Text " Unlike their classical counterparts, confined to a binary fate, qubits waltz through a symphony of probabilities, traversing multiple paths simultaneously" " As civilizations emerged, so did more formalized systems of medicine, such as those practiced in ancient Egypt, Greece, and China"

Replace-Text " This enigmatic ballet of superposition grants quantum computers exponential speed advantages over their classical counterparts, enabling the resolution of previously intractable problems" " These early systems were based on a combination of empirical observations, philosophical beliefs, and religious rituals"

Replace-Text "Beyond the realm of theoretical possibilities, quantum computing is poised to revolutionize medicine, materials science, and artificial intelligence" "Over the centuries, medicine underwent profound changes as new ideas and discoveries emerged"

Replace-Text " Novel drug discoveries, tailored to individual genetic profiles, hold the promise of personalized healthcare, while quantum algorithms illuminate the path towards previously elusive materials with remarkable properties" " The development of the microscope in the 17th century revolutionized our understanding of the human body and disease"

Replace-Text " Artificial intelligence, empowered by quantum computing, embarks on an unprecedented ascent, soaring to new heights of efficiency and accuracy" " The discovery of microorganisms in the 19th century led to the germ theory of disease, which transformed how we approach infection and prevention"

# ---------------------------------------------------------------------
# 6) Body paragraph - brand-new material appended at the end of the
#    paragraph (new sentences + manual line breaks).
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$insPos = $p5.Range.End - 1
$tail = $d.Range($insPos, $insPos)
$tail.InsertAfter("ZZINSERTZZ")

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Replacement.Font.Name = "Times New Roman"
$find.Replacement.Font.Size = 12
$find.Replacement.Font.Color = 0

$newTail = " Advances in chemistry gave rise to new drugs and treatments that saved countless lives." + `
    "^l^lIntroduction Continued:^l^l" + `
    "In the 20th and 21st centuries, medicine has witnessed an explosion of knowledge and innovation." + `
    " The mapping of the human genome, the advent of gene therapy, and the development of personalized medicine have opened up new possibilities for treating and preventing diseases." + `
    " Robotic surgery, minimally invasive techniques, and telemedicine have transformed the way healthcare is delivered.^l" + `
    "Despite these remarkable advancements, the challenges of disease and suffering persist." + `
    " New diseases emerge, and old ones continue to plague humanity." + `
    " The rising burden of chronic conditions presents unique challenges for healthcare systems worldwide." + `
    " The need for accessible, affordable, and equitable healthcare remains a pressing issue.^l^l" + `
    "Introduction Continued:^l^l" + `
    "Medicine's journey is a testament to the resilience of the human spirit and our unwavering commitment to improving lives." + `
    " It is a story of progress and innovation, challenges and triumphs." + `
    " As we continue to unravel the mysteries of the human body, develop new technologies, and embrace the power of collaboration, the future of medicine holds immense promise for a healthier and more compassionate world."

$find.Execute("ZZINSERTZZ", $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2)

# ---------------------------------------------------------------------
# 7) Summary paragraph - sentence-by-sentence rewrite
# ---------------------------------------------------------------------
Replace-Text "The enigmatic tapestry of quantum computing is unraveling before our eyes, revealing a world of limitless potential" "Medicine's history is a tale of human endeavor, innovation, and the pursuit of healing"

Replace-Text " This mind-bending realm, where superposition dances and qubits pirouette in a quantum waltz, promises transformative breakthroughs across diverse fields, from medicine to materials science to artificial intelligence" " From ancient herbal remedies to modern medical marvels, medicine has undergone profound transformations over time, driven by the dedication of individuals committed to alleviating suffering"

Replace-Text " With the dawn of quantum computing, we stand at the precipice of a new era, poised to witness the unfolding of a technological odyssey that will reshape our understanding of reality and redefine the boundaries of human ingenuity" " Despite ongoing challenges, the future of medicine offers hope and promise for a healthier world"

# ---------------------------------------------------------------------
# 8) Trailing empty paragraph added right before the sectPr
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endOfDoc = $lastPara.Range.End
$newParaRange = $d.Range($endOfDoc, $endOfDoc)
$newParaRange.InsertParagraphAfter()

Write-Output "done"
